$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    4  = 1
    5  = 1
    6  = 1
    7  = 0
    8  = 0
    9  = 2
    10 = 2
    11 = 2
    12 = 1
    13 = 2
    14 = 0
    15 = 1
    16 = 1
    17 = 0
    18 = 0
    19 = 2
    20 = 1
    21 = 3
    22 = 1
    23 = 1
    24 = 0
    25 = 0
    26 = 3
    27 = 1
    28 = 0
    29 = 0
    30 = 1
    31 = 1
    32 = 0
    33 = 1
    34 = 0
    35 = 2
    36 = 2
    37 = 2
    38 = 2
    39 = 0
    41 = 1
    42 = 2
    43 = 1
    44 = 1
    45 = 1
    46 = 0
    47 = 0
    48 = 0
    49 = 1
    50 = 0
    51 = 2
    52 = 1
    53 = 1
    54 = 1
    55 = 0
    56 = 3
    57 = 2
    59 = 2
    60 = 1
    61 = 1
    62 = 2
    63 = 2
    64 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
